$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused duplicate "UpdatePassword" column (H) - column G already covers it
$ws.Columns("H").Delete()

# Update header row: C1 stays "Initial" (shared-string re-indexed automatically)

# Update existing row 2 values (email changed) and widen/retype columns D & F
$ws.Range("F2").Value = "abs@gmail.com"

# New data rows 3-12
$ws.Range("A3").Value = "plover@icloud.com"
$ws.Range("B3").Value = "aA1asdfghjk"
$ws.Range("C3").Value = "IN"
$ws.Range("D3").Value = "Ingrid"
$ws.Range("E3").Value = "Nike"
$ws.Range("F3").Value = "pl@gmail.com"
$ws.Range("G3").Value = "Ss2Ghjahsfh"

$ws.Range("A4").Value = "esasaki@icloud.com"
$ws.Range("B4").Value = "aA1asdfghjk"
$ws.Range("C4").Value = "GO"
$ws.Range("D4").Value = "Gogkiud"
$ws.Range("E4").Value = "Oka"
$ws.Range("F4").Value = "es@gomas.com"
$ws.Range("G4").Value = "esS1sjfhasfh"

$ws.Range("A5").Value = "jguyer@msn.com"
$ws.Range("B5").Value = "aA1asdfghjk"
$ws.Range("C5").Value = "OR"
$ws.Range("D5").Value = "Orgrmo"
$ws.Range("E5").Value = "Rmke"
$ws.Range("F5").Value = "jg@gmamg.com"
$ws.Range("G5").Value = "Or1hshafhausf"

$ws.Range("A6").Value = "dieman@live.com"
$ws.Range("B6").Value = "aA1asdfghjk"
$ws.Range("C6").Value = "AS"
$ws.Range("D6").Value = "Asgnlc"
$ws.Range("E6").Value = "Slka"
$ws.Range("F6").Value = "di@gmali.com"
$ws.Range("G6").Value = "Or1hshafhausf"

$ws.Range("A7").Value = "penna@mac.com"
$ws.Range("B7").Value = "aA1asdfghjk"
$ws.Range("C7").Value = "ME"
$ws.Range("D7").Value = "Megmam"
$ws.Range("E7").Value = "Emaka"
$ws.Range("F7").Value = "pe@gmamae.com"
$ws.Range("G7").Value = "esS1sjfhasfh"

$ws.Range("A8").Value = "lahvak@outlook.com"
$ws.Range("B8").Value = "aA1asdfghjk"
$ws.Range("C8").Value = "NO"
$ws.Range("D8").Value = "Nogkook"
$ws.Range("E8").Value = "Ooka"
$ws.Range("F8").Value = "la@gomaoa.com"
$ws.Range("G8").Value = "Aa1qwerty"

$ws.Range("A9").Value = "eabrown@sbcglobal.net"
$ws.Range("B9").Value = "aA1asdfghjk"
$ws.Range("C9").Value = "EP"
$ws.Range("D9").Value = "Epgwnlobal"
$ws.Range("E9").Value = "Pko"
$ws.Range("F9").Value = "ea@gl.netaa.net"
$ws.Range("G9").Value = "Ss2Ghjahsfh"

$ws.Range("A10").Value = "telbij@msn.com"
$ws.Range("B10").Value = "aA1asdfghjk"
$ws.Range("C10").Value = "AT"
$ws.Range("D10").Value = "Atgjmo"
$ws.Range("E10").Value = "Tmki"
$ws.Range("F10").Value = "te@gmame.com"
$ws.Range("G10").Value = "esS1sjfhasfh"

$ws.Range("A11").Value = "yzheng@mac.com"
$ws.Range("B11").Value = "aA1asdfghjk"
$ws.Range("C11").Value = "UP"
$ws.Range("D11").Value = "Upggmo"
$ws.Range("E11").Value = "Pmkn"
$ws.Range("F11").Value = "yz@gmamz.com"
$ws.Range("G11").Value = "Or1hshafhausf"

$ws.Range("A12").Value = "benits@sbcglobal.net"
$ws.Range("B12").Value = "aA1asdfghjk"
$ws.Range("C12").Value = "AM"
$ws.Range("D12").Value = "Amgssobal"
$ws.Range("E12").Value = "Mskt"
$ws.Range("F12").Value = "be@g.netase.net"
$ws.Range("G12").Value = "esS1sjfhasfh"

# Rebuild hyperlinks from scratch (collection delete clears all, then re-add in the
# same order as the target file: A2, F2, F3..F12)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:austinryang1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:abs@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:pl@gmail.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:es@gomas.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:jg@gmamg.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:di@gmali.com")
$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:pe@gmamae.com")
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:la@gomaoa.com")
$ws.Hyperlinks.Add($ws.Range("F9"), "mailto:ea@gl.netaa.net")
$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:te@gmame.com")
$ws.Hyperlinks.Add($ws.Range("F11"), "mailto:yz@gmamz.com")
$ws.Hyperlinks.Add($ws.Range("F12"), "mailto:be@g.netase.net")

# Restore the plain "Hyperlink" cell style on every linked cell (Hyperlinks.Add
# leaves a stray extra style otherwise)
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F12").Style = "Hyperlink"

# Column width / formatting tweaks
$ws.Columns("D").ColumnWidth = 9.3
$ws.Columns("F").ColumnWidth = 17.8

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection as left by the editing session
$ws.Range("J12").Select()
